$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.678.23'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '3.451.04'
$ws.Range('D5').Value = '592.02'
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('D6').Value = '178.86'
$ws.Range('E6').Value = '  -2.73%  '
$ws.Range('E7').Value = '  +1.94%  '
$ws.Range('D9').Value = '3.448.67'
$ws.Range('E9').Value = '  -2.32%  '
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('E12').Value = '  -3.80%  '
$ws.Range('D13').Value = '4.045.38'
$ws.Range('E13').Value = '  -2.29%  '
$ws.Range('D14').Value = '32.04'
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('E15').Value = '  -0.73%  '
$ws.Range('D16').Value = '67.678.68'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('E17').Value = '  -3.27%  '
$ws.Range('D18').Value = '3.450.66'
$ws.Range('E18').Value = '  -2.07%  '
$ws.Range('E19').Value = '  -4.24%  '
$ws.Range('D20').Value = '14.00'
$ws.Range('E20').Value = '  -6.38%  '
$ws.Range('D21').Value = '388.94'
$ws.Range('D22').Value = '7.85'
$ws.Range('E22').Value = '  -3.40%  '
$ws.Range('D23').Value = '5.82'
$ws.Range('E23').Value = '  +1.96%  '
$ws.Range('D24').Value = '0.996'
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('E25').Value = '  -2.53%  '
$ws.Range('D26').Value = '71.44'
$ws.Range('E26').Value = '  -3.18%  '
$ws.Range('E27').Value = '  -4.87%  '
$ws.Range('D28').Value = '10.24'
$ws.Range('E28').Value = '  -4.86%  '
$ws.Range('E29').Value = '  -2.38%  '
$ws.Range('D30').Value = '0.997'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').Value = '6.06'
$ws.Range('E31').Value = '  -3.58%  '
$ws.Range('E32').Value = '  -1.90%  '
$ws.Range('E33').Value = '  -5.50%  '
$ws.Range('D34').Value = '23.37'
$ws.Range('E34').Value = '  -3.42%  '
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').Value = '7.22'
$ws.Range('E36').Value = '  -3.29%  '
$ws.Range('E37').Value = '  -7.67%  '
$ws.Range('D38').Value = '160.50'
$ws.Range('E38').Value = '  -2.14%  '
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('E40').Value = '  -5.19%  '
$ws.Range('D41').Value = '2.74'
$ws.Range('E41').Value = '  -2.08%  '
$ws.Range('D42').Value = '4.60'
$ws.Range('E42').Value = '  -3.80%  '
$ws.Range('E43').Value = '  -7.80%  '
$ws.Range('D44').Value = '25.88'
$ws.Range('E44').Value = '  -4.95%  '
$ws.Range('E45').Value = '  -3.99%  '
$ws.Range('E46').Value = '  -6.05%  '
$ws.Range('D47').Value = '2.699.31'
$ws.Range('E47').Value = '  -6.63%  '
$ws.Range('D48').Value = '41.19'
$ws.Range('E48').Value = '  -3.29%  '
$ws.Range('E49').Value = '  -3.52%  '
$ws.Range('D50').Value = '323.41'
$ws.Range('E50').Value = '  -8.57%  '
$ws.Range('E51').Value = '  -4.81%  '
